$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.854.78'
$ws.Range('E2').Value = '  -6.75%  '
$ws.Range('D3').Value = '3.694.55'
$ws.Range('E3').Value = '  -6.27%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '579.45'
$ws.Range('E5').Value = '  -3.82%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '172.04'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('D7').Value = '3.686.59'
$ws.Range('E7').Value = '  -6.38%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.624'
$ws.Range('E8').Value = '  -8.85%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.996'
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.703'
$ws.Range('E10').Value = '  -10.74%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.161'
$ws.Range('E11').Value = '  -13.13%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '51.01'
$ws.Range('E12').Value = '  -9.22%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000288'
$ws.Range('E13').Value = '  -13.34%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.36'
$ws.Range('E14').Value = '  -10.40%  '
$ws.Range('D15').Value = '4.270.68'
$ws.Range('E15').Value = '  -6.35%  '
$ws.Range('D16').Value = '3.690.39'
$ws.Range('E16').Value = '  -6.43%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '19.31'
$ws.Range('E17').Value = '  -10.40%  '
$ws.Range('E18').Value = '  -3.28%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.77'
$ws.Range('E19').Value = '  -9.81%  '
$ws.Range('E20').Value = '  -9.66%  '
$ws.Range('D21').Value = '67.534.54'
$ws.Range('E21').Value = '  -6.96%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '404.63'
$ws.Range('E22').Value = '  -8.82%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.47'
$ws.Range('E23').Value = '  -5.66%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '87.91'
$ws.Range('E24').Value = '  -8.12%  '
$ws.Range('E25').Value = '  -9.50%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.67'
$ws.Range('E26').Value = '  -10.09%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.76'
$ws.Range('E27').Value = '  -3.46%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.02'
$ws.Range('E28').Value = '  +1.28%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '3.78'
$ws.Range('E29').Value = '  -12.18%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.44'
$ws.Range('E30').Value = '  -9.47%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '32.46'
$ws.Range('E31').Value = '  -9.49%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.39'
$ws.Range('E32').Value = '  -5.85%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '12.37'
$ws.Range('E33').Value = '  -10.67%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '64.73'
$ws.Range('E34').Value = '  -5.60%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.115'
$ws.Range('E35').Value = '  -9.84%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '603.39'
$ws.Range('E36').Value = '  -4.66%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '42.96'
$ws.Range('E37').Value = '  -15.36%  '
$ws.Range('D38').Value = '0.0₃0884'
$ws.Range('E38').Value = '  -11.93%  '
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.393'
$ws.Range('E40').Value = '  -7.89%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.136'
$ws.Range('E42').Value = '  -6.82%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.76'
$ws.Range('E43').Value = '  +5.81%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.97'
$ws.Range('E44').Value = '  -12.22%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0433'
$ws.Range('E45').Value = '  -9.44%  '
$ws.Range('E46').Value = '  -13.47%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.18'
$ws.Range('E47').Value = '  -13.05%  '
$ws.Range('D48').Value = '2.809.21'
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.133'
$ws.Range('E49').Value = '  -9.71%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.68'
$ws.Range('E50').Value = '  -5.86%  '
$ws.Range('E51').Value = '  -7.99%  '
